$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values to the "custom accuracy" (2 decimal place) figures ---
$ws.Cells.Item(5, 2).Value  = 17.17               # B5
$ws.Cells.Item(5, 3).Value  = 12.83               # C5
$ws.Cells.Item(5, 4).Value  = 1.09                # D5
$ws.Cells.Item(5, 5).Value  = 37.63               # E5
$ws.Cells.Item(5, 6).Value  = 30.6                # F5
$ws.Cells.Item(5, 7).Value  = 13.18               # G5
$ws.Cells.Item(5, 8).Value  = 52.04               # H5
$ws.Cells.Item(5, 9).Value  = 20.88               # I5
$ws.Cells.Item(5, 10).Value = 9.449999999999999   # J5
$ws.Cells.Item(5, 11).Value = 13.52               # K5
$ws.Cells.Item(5, 12).Value = 15.08               # L5
$ws.Cells.Item(5, 13).Value = 16.09               # M5
$ws.Cells.Item(5, 14).Value = 4.51                # N5
$ws.Cells.Item(5, 15).Value = 13.53               # O5
$ws.Cells.Item(5, 16).Value = 19.14               # P5
$ws.Cells.Item(5, 17).Value = 11.52               # Q5
$ws.Cells.Item(5, 18).Value = 0.47                # R5
$ws.Cells.Item(5, 19).Value = 0.66                # S5
$ws.Cells.Item(5, 20).Value = 198.69              # T5
$ws.Cells.Item(5, 21).Value = 37.84               # U5
$ws.Cells.Item(5, 22).Value = 12.49               # V5
$ws.Cells.Item(5, 23).Value = 25.33               # W5
$ws.Cells.Item(5, 24).Value = 13.35               # X5
$ws.Cells.Item(5, 25).Value = 1.76                # Y5
$ws.Cells.Item(5, 26).Value = 25.98               # Z5
$ws.Cells.Item(5, 27).Value = 11.03               # AA5
$ws.Cells.Item(5, 28).Value = 9.81                # AB5
$ws.Cells.Item(5, 29).Value = 11.51               # AC5
$ws.Cells.Item(5, 30).Value = 15.9                # AD5
$ws.Cells.Item(5, 31).Value = 0.48                # AE5
$ws.Cells.Item(5, 32).Value = 47.36               # AF5
$ws.Cells.Item(5, 33).Value = 6.98                # AG5
$ws.Cells.Item(5, 34).Value = 15.62               # AH5

# --- Remove row 6 entirely (data trimmed to fewer sample rows) ---
$ws.Rows(6).Delete()

# --- Narrow a handful of columns from width 8 to width 7 ---
# (COM ColumnWidth is offset from the raw OOXML character width by ~0.83
#  for the default Calibri 11 font, so request 6.17 to land on raw width 7)
$ws.Columns(3).ColumnWidth = 6.17    # C
$ws.Columns(17).ColumnWidth = 6.17   # Q
$ws.Columns(22).ColumnWidth = 6.17   # V
$ws.Columns(27).ColumnWidth = 6.17   # AA
$ws.Columns(29).ColumnWidth = 6.17   # AC
